$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.06729733333333333
$ws.Range("H2").Value = 0.201892
$ws.Range("I2").Value = 0.01373511018321553
$ws.Range("J2").Value = 0.01373511018321553
$ws.Range("M2").Value = 13.17295566666667
$ws.Range("N2").Value = 39.518867
$ws.Range("O2").Value = 0.133784132206724
$ws.Range("P2").Value = 0.133784132206724
$ws.Range("Q2").Value = 0.8865047884848889
$ws.Range("R2").Value = 7.978543096364
$ws.Range("S2").Value = 0.001837539796625228
$ws.Range("T2").Value = 0.001837539796625228
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.06729733333333333
$ws.Range("H3").Value = 0.201892
$ws.Range("I3").Value = 0.01373511018321553
$ws.Range("J3").Value = 0.01373511018321553
$ws.Range("O3").Value = 0.4382627974978752
$ws.Range("P3").Value = 0.4382627974978752
$ws.Range("Q3").Value = 2.904096787773778
$ws.Range("R3").Value = 26.136871089964
$ws.Range("S3").Value = 0.006019587812837592
$ws.Range("T3").Value = 0.006019587812837591
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.06729733333333333
$ws.Range("H4").Value = 0.201892
$ws.Range("I4").Value = 0.01373511018321553
$ws.Range("J4").Value = 0.01373511018321553
$ws.Range("M4").Value = 21.06166566666667
$ws.Range("N4").Value = 63.184997
$ws.Range("O4").Value = 0.2139016281041017
$ws.Range("P4").Value = 0.2139016281041017
$ws.Range("Q4").Value = 1.417393934924889
$ws.Range("R4").Value = 12.756545414324
$ws.Range("S4").Value = 0.002937962430379029
$ws.Range("T4").Value = 0.002937962430379029
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.06729733333333333
$ws.Range("H5").Value = 0.201892
$ws.Range("I5").Value = 0.01373511018321553
$ws.Range("J5").Value = 0.01373511018321553
$ws.Range("M5").Value = 21.076417
$ws.Range("N5").Value = 63.229251
$ws.Range("O5").Value = 0.214051442191299
$ws.Range("P5").Value = 0.214051442191299
$ws.Range("Q5").Value = 1.418386660321334
$ws.Range("R5").Value = 12.765479942892
$ws.Range("S5").Value = 0.002940020143373683
$ws.Range("T5").Value = 0.002940020143373682
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 4.83236
$ws.Range("H6").Value = 14.49708
$ws.Range("I6").Value = 0.9862648898167845
$ws.Range("J6").Value = 0.9862648898167844
$ws.Range("M6").Value = 13.17295566666667
$ws.Range("N6").Value = 39.518867
$ws.Range("O6").Value = 0.133784132206724
$ws.Range("P6").Value = 0.133784132206724
$ws.Range("Q6").Value = 63.65646404537334
$ws.Range("R6").Value = 572.90817640836
$ws.Range("S6").Value = 0.1319465924100988
$ws.Range("T6").Value = 0.1319465924100988
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 4.83236
$ws.Range("H7").Value = 14.49708
$ws.Range("I7").Value = 0.9862648898167845
$ws.Range("J7").Value = 0.9862648898167844
$ws.Range("O7").Value = 0.4382627974978752
$ws.Range("P7").Value = 0.4382627974978752
$ws.Range("Q7").Value = 208.5319054747067
$ws.Range("R7").Value = 1876.78714927236
$ws.Range("S7").Value = 0.4322432096850376
$ws.Range("T7").Value = 0.4322432096850375
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 4.83236
$ws.Range("H8").Value = 14.49708
$ws.Range("I8").Value = 0.9862648898167845
$ws.Range("J8").Value = 0.9862648898167844
$ws.Range("M8").Value = 21.06166566666667
$ws.Range("N8").Value = 63.184997
$ws.Range("O8").Value = 0.2139016281041017
$ws.Range("P8").Value = 0.2139016281041017
$ws.Range("Q8").Value = 101.7775507009733
$ws.Range("R8").Value = 915.9979563087601
$ws.Range("S8").Value = 0.2109636656737227
$ws.Range("T8").Value = 0.2109636656737227
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 4.83236
$ws.Range("H9").Value = 14.49708
$ws.Range("I9").Value = 0.9862648898167845
$ws.Range("J9").Value = 0.9862648898167844
$ws.Range("M9").Value = 21.076417
$ws.Range("N9").Value = 63.229251
$ws.Range("O9").Value = 0.214051442191299
$ws.Range("P9").Value = 0.214051442191299
$ws.Range("Q9").Value = 101.84883445412
$ws.Range("R9").Value = 916.6395100870801
$ws.Range("S9").Value = 0.2111114220479254
$ws.Range("T9").Value = 0.2111114220479254
